$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "Erkek baggy pantolon, bol ve rahat kesimiyle öne çıkan, modern sokak stilinin vazgeçilmez parçasıdır. Kalçadan paçaya kadar geniş formu sayesinde hareket özgürlüğü sunar.Ürün içeriği 100% Pamuk. Rahatlığından ödün vermeden tarzını yansıtmak isteyen erkekler için ideal bir seçimdir.29-38 Beden beden seçeneği mevcuttur."

# Add a new product row (Baggy Füme) first
$ws.Range("A5").Value = "Baggy Füme"
$ws.Range("B5").Value = "300 tl"
$ws.Range("C5").Value = "Jeans"
$ws.Range("D5").Value = "BAG3.jpg"
$ws.Range("E5").Value = $desc
$ws.Range("F5").Value = "Var"

# Update the "aciklama" (description) column for the existing three rows
# from the short "100% Pamuk" text to the full product description.
$ws.Range("E2").Value = $desc
$ws.Range("E3").Value = $desc
$ws.Range("E4").Value = $desc

# Give column E an explicit width (matches the new customWidth seen for col 5)
$ws.Columns.Item(5).ColumnWidth = 8.3

# Move the selection, matching the last active cell recorded in the file
$ws.Range("F14").Select()
